$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Overall"
$ws.Cells.Item(2, 2).Value = 1212
$ws.Cells.Item(2, 3).Value = 4864.78300330033
$ws.Cells.Item(2, 4).Value = 4030.62673192859
$ws.Cells.Item(2, 5).Value = 5698.93927467208
$ws.Cells.Item(2, 6).Value = 64.8884217351934
$ws.Cells.Item(2, 7).Value = 53.9111984327623
$ws.Cells.Item(2, 8).Value = 76.6485603333175
$ws.Cells.Item(3, 1).Value = "Acinobacter"
$ws.Cells.Item(3, 2).Value = 93
$ws.Cells.Item(3, 3).Value = 6124.32258064516
$ws.Cells.Item(3, 4).Value = 3226.79203284396
$ws.Cells.Item(3, 5).Value = 9021.85312844636
$ws.Cells.Item(3, 6).Value = 64.6167943642757
$ws.Cells.Item(3, 7).Value = 36.1918016499488
$ws.Cells.Item(3, 8).Value = 98.9744511671965
$ws.Cells.Item(4, 1).Value = "klebsiella"
$ws.Cells.Item(4, 2).Value = 207
$ws.Cells.Item(4, 3).Value = 7897.14975845411
$ws.Cells.Item(4, 4).Value = 5602.10540821417
$ws.Cells.Item(4, 5).Value = 10192.194108694
$ws.Cells.Item(4, 6).Value = 96.8953873520834
$ws.Cells.Item(4, 7).Value = 69.7943364112936
$ws.Cells.Item(4, 8).Value = 128.322065269713
$ws.Cells.Item(5, 1).Value = "Clostridium"
$ws.Cells.Item(5, 2).Value = 27
$ws.Cells.Item(5, 3).Value = 2461.88888888889
$ws.Cells.Item(5, 4).Value = -198.034289005532
$ws.Cells.Item(5, 5).Value = 5121.81206678331
$ws.Cells.Item(5, 6).Value = 72.3972877312904
$ws.Cells.Item(5, 7).Value = 15.4373426626473
$ws.Cells.Item(5, 8).Value = 157.462829025449
$ws.Cells.Item(6, 1).Value = "Enterococcus"
$ws.Cells.Item(6, 2).Value = 441
$ws.Cells.Item(6, 3).Value = 5090.07936507936
$ws.Cells.Item(6, 4).Value = 3754.05071928177
$ws.Cells.Item(6, 5).Value = 6426.10801087696
$ws.Cells.Item(6, 6).Value = 65.5078750545022
$ws.Cells.Item(6, 7).Value = 48.8658424320023
$ws.Cells.Item(6, 8).Value = 84.010356288206
$ws.Cells.Item(7, 1).Value = "Escherichiacoli"
$ws.Cells.Item(7, 2).Value = 247
$ws.Cells.Item(7, 3).Value = 4285.68825910931
$ws.Cells.Item(7, 4).Value = 2472.06516561117
$ws.Cells.Item(7, 5).Value = 6099.31135260746
$ws.Cells.Item(7, 6).Value = 51.2030403669678
$ws.Cells.Item(7, 7).Value = 31.7394759675061
$ws.Cells.Item(7, 8).Value = 73.5422070591352
$ws.Cells.Item(8, 1).Value = "Pseudomonas"
$ws.Cells.Item(8, 2).Value = 140
$ws.Cells.Item(8, 3).Value = 5383.36428571429
$ws.Cells.Item(8, 4).Value = 2827.14649923311
$ws.Cells.Item(8, 5).Value = 7939.58207219546
$ws.Cells.Item(8, 6).Value = 68.3921251544791
$ws.Cells.Item(8, 7).Value = 38.6970089929833
$ws.Cells.Item(8, 8).Value = 104.444984213584
$ws.Cells.Item(9, 1).Value = "Candida"
$ws.Cells.Item(9, 2).Value = 392
$ws.Cells.Item(9, 3).Value = 4379.93112244898
$ws.Cells.Item(9, 4).Value = 3148.56358469062
$ws.Cells.Item(9, 5).Value = 5611.29866020734
$ws.Cells.Item(9, 6).Value = 67.368130068063
$ws.Cells.Item(9, 7).Value = 50.3133693989219
$ws.Cells.Item(9, 8).Value = 86.3579472304809
$ws.Cells.Item(10, 1).Value = "Staphylococcus"
$ws.Cells.Item(10, 2).Value = 454
$ws.Cells.Item(10, 3).Value = 7023.04185022026
$ws.Cells.Item(10, 4).Value = 5646.598832422
$ws.Cells.Item(10, 5).Value = 8399.48486801854
$ws.Cells.Item(10, 6).Value = 87.2034416652225
$ws.Cells.Item(10, 7).Value = 69.9099582948445
$ws.Cells.Item(10, 8).Value = 106.257060639675
$ws.Cells.Item(11, 1).Value = "Blood"
$ws.Cells.Item(11, 2).Value = 417
$ws.Cells.Item(11, 3).Value = 9736.36211031175
$ws.Cells.Item(11, 4).Value = 8136.12689722172
$ws.Cells.Item(11, 5).Value = 11336.5973234018
$ws.Cells.Item(11, 6).Value = 127.167474230522
$ws.Cells.Item(11, 7).Value = 106.403026976817
$ws.Cells.Item(11, 8).Value = 150.020855334022
$ws.Cells.Item(12, 1).Value = "Urinary"
$ws.Cells.Item(12, 2).Value = 592
$ws.Cells.Item(12, 3).Value = 3477.55236486487
$ws.Cells.Item(12, 4).Value = 2587.63700177929
$ws.Cells.Item(12, 5).Value = 4367.46772795044
$ws.Cells.Item(12, 6).Value = 51.2598517950677
$ws.Cells.Item(12, 7).Value = 39.4830567718202
$ws.Cells.Item(12, 8).Value = 64.0309819313352
$ws.Cells.Item(13, 1).Value = "Respiratory"
$ws.Cells.Item(13, 2).Value = 267
$ws.Cells.Item(13, 3).Value = 8593.22846441948
$ws.Cells.Item(13, 4).Value = 6870.42788015657
$ws.Cells.Item(13, 5).Value = 10316.0290486824
$ws.Cells.Item(13, 6).Value = 95.349073603191
$ws.Cells.Item(13, 7).Value = 74.6299266386883
$ws.Cells.Item(13, 8).Value = 118.526465034719
$ws.Cells.Item(14, 1).Value = "Wound"
$ws.Cells.Item(14, 2).Value = 170
$ws.Cells.Item(14, 3).Value = 7376.07647058824
$ws.Cells.Item(14, 4).Value = 5262.32772055414
$ws.Cells.Item(14, 5).Value = 9489.82522062233
$ws.Cells.Item(14, 6).Value = 110.185676953546
$ws.Cells.Item(14, 7).Value = 78.4682011470801
$ws.Cells.Item(14, 8).Value = 147.540001593967
$ws.Cells.Item(15, 1).Value = "Rectal"
$ws.Cells.Item(15, 2).Value = 62
$ws.Cells.Item(15, 3).Value = 8144.58064516129
$ws.Cells.Item(15, 4).Value = 3732.43330338389
$ws.Cells.Item(15, 5).Value = 12556.7279869387
$ws.Cells.Item(15, 6).Value = 133.431450505943
$ws.Cells.Item(15, 7).Value = 0.458296605095265
$ws.Cells.Item(15, 8).Value = 244.574433398863

Write-Output "Applied updated risultati data"
